$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 63803.6
$ws.Range("J7").Value = 63803.6
$ws.Range("L7").Value = 63803.6
$ws.Range("N7").Value = -64027.6
$ws.Range("H14").Value = 63803.6
$ws.Range("J14").Value = 63803.6
$ws.Range("L14").Value = 63803.6
$ws.Range("N14").Value = -64185.6
$ws.Range("H15").Value = 1290.6171
$ws.Range("I15").Value = 1290.6171
$ws.Range("K15").Value = 3871.8513
$ws.Range("M15").Value = -3702.8513
$ws.Range("H87").Value = 25882
$ws.Range("J87").Value = 25882
$ws.Range("L87").Value = 25882
$ws.Range("N87").Value = -28378
$ws.Range("H90").Value = 25882
$ws.Range("J90").Value = 25882
$ws.Range("L90").Value = 77646
$ws.Range("N90").Value = -90126
$ws.Range("H137").Value = 5232.4116
$ws.Range("I137").Value = 3264.0667
$ws.Range("J137").Value = 19995
$ws.Range("K137").Value = 9792.2001
$ws.Range("L137").Value = 59985
$ws.Range("M137").Value = -7242.2001
$ws.Range("N137").Value = -65085
$ws.Range("H138").Value = 2272.9412
$ws.Range("J138").Value = 2518.7144
$ws.Range("L138").Value = 7556.1432
$ws.Range("N138").Value = -17836.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2720.5
$ws.Range("I2").Value = 2623
$ws.Range("J2").Value = 3013
$ws.Range("K2").Value = 2623
$ws.Range("L2").Value = 3013
$ws.Range("M2").Value = -2510
$ws.Range("N2").Value = -3239
$ws.Range("H61").Value = 3736.05
$ws.Range("I61").Value = 2946.7778
$ws.Range("K61").Value = 2946.7778
$ws.Range("M61").Value = -2734.7778
$ws.Range("H74").Value = 1405.3125
$ws.Range("I74").Value = 1052.7693
$ws.Range("K74").Value = 1052.7693
$ws.Range("M74").Value = -178.7692999999999
$ws.Range("H77").Value = 1405.3125
$ws.Range("I77").Value = 1052.7693
$ws.Range("K77").Value = 5263.8465
$ws.Range("M77").Value = -895.8464999999997
$ws.Range("H116").Value = 2720.5
$ws.Range("I116").Value = 2623
$ws.Range("J116").Value = 3013
$ws.Range("K116").Value = 2623
$ws.Range("L116").Value = 3013
$ws.Range("M116").Value = -329
$ws.Range("N116").Value = -7601
$ws.Range("H132").Value = 3832.6858
$ws.Range("I132").Value = 2463.8333
$ws.Range("K132").Value = 7391.499899999999
$ws.Range("M132").Value = -4861.499899999999
$ws.Range("H136").Value = 3736.05
$ws.Range("I136").Value = 2946.7778
$ws.Range("K136").Value = 8840.3334
$ws.Range("M136").Value = -6290.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2720.5
$ws.Range("I3").Value = 2623
$ws.Range("J3").Value = 3013
$ws.Range("K3").Value = 2623
$ws.Range("L3").Value = 3013
$ws.Range("M3").Value = -2509
$ws.Range("N3").Value = -3241
$ws.Range("H134").Value = 3194.1428
$ws.Range("I134").Value = 3152.3809
$ws.Range("J134").Value = 3319.4285
$ws.Range("K134").Value = 9457.1427
$ws.Range("L134").Value = 9958.2855
$ws.Range("M134").Value = -6922.1427
$ws.Range("N134").Value = -15028.2855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 28091.092
$ws.Range("J4").Value = 28091.092
$ws.Range("L4").Value = 28091.092
$ws.Range("N4").Value = -28315.092
$ws.Range("H31").Value = 5749.2383
$ws.Range("I31").Value = 1039.8286
$ws.Range("J31").Value = 11636
$ws.Range("K31").Value = 1039.8286
$ws.Range("L31").Value = 11636
$ws.Range("M31").Value = -744.8286000000001
$ws.Range("N31").Value = -12226
$ws.Range("H34").Value = 5749.2383
$ws.Range("I34").Value = 1039.8286
$ws.Range("J34").Value = 11636
$ws.Range("K34").Value = 1039.8286
$ws.Range("L34").Value = 11636
$ws.Range("M34").Value = -837.8286000000001
$ws.Range("N34").Value = -12040
$ws.Range("H58").Value = 1676.8667
$ws.Range("I58").Value = 1436.5834
$ws.Range("J58").Value = 2638
$ws.Range("K58").Value = 1436.5834
$ws.Range("L58").Value = 2638
$ws.Range("M58").Value = -1233.5834
$ws.Range("N58").Value = -3044
$ws.Range("H99").Value = 1842.1052
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H126").Value = 1842.1052
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 6946514
$ws.Range("I132").Value = 1758.55
$ws.Range("J132").Value = 41670292
$ws.Range("K132").Value = 5275.65
$ws.Range("L132").Value = 125010876
$ws.Range("M132").Value = -2745.65
$ws.Range("N132").Value = -125015936
$ws.Range("H134").Value = 4393.0293
$ws.Range("I134").Value = 4442.3706
$ws.Range("K134").Value = 13327.1118
$ws.Range("M134").Value = -10792.1118
$ws.Range("H136").Value = 1676.8667
$ws.Range("I136").Value = 1436.5834
$ws.Range("J136").Value = 2638
$ws.Range("K136").Value = 4309.7502
$ws.Range("L136").Value = 7914
$ws.Range("M136").Value = -1759.7502
$ws.Range("N136").Value = -13014

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 833.975
$ws.Range("I113").Value = 586.4737
$ws.Range("J113").Value = 1057.9048
$ws.Range("K113").Value = 1759.4211
$ws.Range("L113").Value = 3173.7144
$ws.Range("M113").Value = 410.5789
$ws.Range("N113").Value = -7513.7144
$ws.Range("H137").Value = 7584667.5
$ws.Range("I137").Value = 20848784
$ws.Range("J137").Value = 5172.2144
$ws.Range("K137").Value = 62546352
$ws.Range("L137").Value = 15516.6432
$ws.Range("M137").Value = -62541252
$ws.Range("N137").Value = -25716.6432
$ws.Range("H140").Value = 1803.5454
$ws.Range("I140").Value = 1057.8
$ws.Range("K140").Value = 3173.4
$ws.Range("M140").Value = 2006.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H132").Value = 3500.5715
$ws.Range("I132").Value = 3300.9333
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 9902.7999
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -7372.7999
$ws.Range("N132").Value = -17059.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1713334
$ws.Range("J2").Value = 1713334
$ws.Range("L2").Value = 1713334
$ws.Range("N2").Value = -1713558
$ws.Range("H132").Value = 3489.2104
$ws.Range("I132").Value = 2820.0667
$ws.Range("K132").Value = 8460.2001
$ws.Range("M132").Value = -5930.2001
$ws.Range("H136").Value = 12822852
$ws.Range("I136").Value = 2399.6365
$ws.Range("J136").Value = 83335336
$ws.Range("K136").Value = 7198.9095
$ws.Range("L136").Value = 250006008
$ws.Range("M136").Value = -4648.9095
$ws.Range("N136").Value = -250011108

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 204028000
$ws.Range("J2").Value = 333380000
$ws.Range("L2").Value = 333380000
$ws.Range("N2").Value = -333380224
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
$ws.Range("H126").Value = 2111.4443
$ws.Range("I126").Value = 680.6
$ws.Range("J126").Value = 3900
$ws.Range("K126").Value = 2041.8
$ws.Range("L126").Value = 11700
$ws.Range("M126").Value = 428.1999999999998
$ws.Range("N126").Value = -16640
$ws.Range("H132").Value = 4904326.5
$ws.Range("I132").Value = 2684.2632
$ws.Range("J132").Value = 11113073
$ws.Range("K132").Value = 8052.7896
$ws.Range("L132").Value = 33339219
$ws.Range("M132").Value = -5522.7896
$ws.Range("N132").Value = -33344279
$ws.Range("H136").Value = 2255.96
$ws.Range("I136").Value = 1925.3611
$ws.Range("J136").Value = 3106.0715
$ws.Range("K136").Value = 5776.0833
$ws.Range("L136").Value = 9318.2145
$ws.Range("M136").Value = -3226.0833
$ws.Range("N136").Value = -14418.2145
